$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume/name/link data to the latest snapshot.
# Column D (Price) values that look numeric must be forced to remain text,
# matching the source data which stores prices as plain strings (e.g. "143.50"
# must not become the number 143.5). We set the number format to Text before
# assigning the value, then restore the default "Normal" style so no stray
# cell formatting is left behind.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.130.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.95%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.482.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.89%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.482.77"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("E9").Value = "  -0.57%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "8.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.89%  "

# Row 11
$ws.Range("E11").Value = "  -5.09%  "

# Row 12
$ws.Range("E12").Value = "  -2.96%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.070.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "

# Row 14
$ws.Range("E14").Value = "  -3.23%  "

# Row 15
$ws.Range("E15").Value = "  -5.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.480.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.160.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.96%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.41%  "

# Row 20
$ws.Range("E20").Value = "  -4.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.38%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.16%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.594"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.30%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.42%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000115"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.96%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.82%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.25%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.26%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.165"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.81%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.26%  "

# Row 34
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.476.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.47%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.51%  "

# Row 38
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.86%  "

# Row 39
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "171.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.67%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0857"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.74%  "

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.78%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.879"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.76%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "45.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.25%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.93%  "

# Row 47
$ws.Range("E47").Value = "  -1.49%  "

# Row 48
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.67%  "

# Row 50
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.933"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.66%  "

# Row 51
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.236"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.71%  "

